$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-06-17 Tuesday"; new = "2025-06-18 Wednesday"},
    @{old = "45×32="; new = "57×44="},
    @{old = "28×81="; new = "57×51="},
    @{old = "77×64="; new = "49×24="},
    @{old = "44×49="; new = "99×75="},
    @{old = "21×22="; new = "33×64="},
    @{old = "14×76="; new = "66×19="},
    @{old = "63×73="; new = "46×79="},
    @{old = "45×87="; new = "11×25="},
    @{old = "85×90="; new = "92×98="},
    @{old = "38×54="; new = "98×75="},
    @{old = "43×83="; new = "61×63="},
    @{old = "28×38="; new = "46×14="},
    @{old = "14×33="; new = "39×73="},
    @{old = "81×40="; new = "36×36="},
    @{old = "55×41="; new = "94×32="},
    @{old = "49×20="; new = "53×91="},
    @{old = "19×15="; new = "22×14="},
    @{old = "56×22="; new = "69×51="},
    @{old = "77×71="; new = "47×51="},
    @{old = "30×72="; new = "25×69="},
    @{old = "23×36="; new = "78×81="},
    @{old = "71×56="; new = "73×39="},
    @{old = "84×15="; new = "33×21="},
    @{old = "73×60="; new = "75×42="},
    @{old = "72×21="; new = "56×53="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
